# Fix the placeholder git-clone URL in the install instructions with the
# real workshop repository URL.
$d = $word.ActiveDocument

$find = "[workshop-git-URL].git"
$replace = "https://github.com/ocean-tracking-network/2023-canssi-ecr-workshop.git"

$d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
